# reviewdb.xlsx — row 16's "time" column (E16) was stored as a raw numeric
# Excel date serial (43592.4854166667) with its own custom date number
# format (style index 3 / numFmtId 165). Every other row in that column
# stores the review time as plain text (e.g. "27/5/2019 15:57" in E2,
# "27/5/2019 15:59" already reused by several other rows). Bring E16 in
# line with the rest of the column: plain text "27/5/2019 15:59" using the
# default (General) cell style instead of the bespoke date format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from a neighboring cell that already uses the plain
# "General" style shared by the rest of column E (e.g. E15), so E16 stops
# referencing the now-unused custom date number format.
$ws.Range("E15").Copy()
$ws.Range("E16").PasteSpecial(-4122)  # xlPasteFormats

# Write the review time as text, matching the other rows in the column.
$ws.Range("E16").Value = "27/5/2019 15:59"
